$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SEGURIDAD")
$b3 = $ws.Range("B3")
$g3 = $ws.Range("G3")
$g3.Value = "TOTAL"
$b3.Copy()
$g3.PasteSpecial(-4122)
$g3.Borders.Item(8).LineStyle = -4142
